$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6, 7, 9, 10, 11 get their data reshuffled (row 8 is untouched).
# New row 6  <- old row 9
# New row 7  <- old row 10
# New row 9  <- old row 11
# New row 10 <- old row 7
# New row 11 <- old row 6
$rows = @(6, 7, 9, 10, 11)
$mapping = @{
    6  = 9
    7  = 10
    9  = 11
    10 = 7
    11 = 6
}

# Capture the current (pre-edit) values first so the writes below don't
# clobber data we still need to read for a later row.
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{
        A = $ws.Range("A$r").Value()
        B = $ws.Range("B$r").Value()
        D = $ws.Range("D$r").Value()
        E = $ws.Range("E$r").Value()
        F = $ws.Range("F$r").Value()
        G = $ws.Range("G$r").Value()
        H = $ws.Range("H$r").Value()
        Q = $ws.Range("Q$r").Value()
        R = $ws.Range("R$r").Value()
    }
}

# The (normally empty/unused) "Bestämningsmetod" cell in column AF travels
# along with the row that owns it: it currently sits on row 7 and needs to
# end up on row 10 (the row that inherits row 7's other data).
$ws.Range("AF7").Copy($ws.Range("AF10"))
$ws.Range("AF7").ClearContents()

foreach ($dest in $rows) {
    $src = $mapping[$dest]
    $vals = $orig[$src]

    $ws.Range("A$dest").Value = $vals.A
    $ws.Range("B$dest").Value = $vals.B
    $ws.Range("D$dest").Value = $vals.D
    $ws.Range("E$dest").Value = $vals.E
    $ws.Range("F$dest").Value = $vals.F
    $ws.Range("G$dest").Value = $vals.G
    $ws.Range("H$dest").Value = $vals.H
    $ws.Range("Q$dest").Value = $vals.Q
    $ws.Range("R$dest").Value = $vals.R
}
